$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.145.28'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.42%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.869.64'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.93%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.40'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.53%  '

$ws.Range("E6").Value = '  -0.19%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5009'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.25%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3897'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.44%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.09639'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +25.29%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.139'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.16%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '40.97'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.60%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.459'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.33%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.90'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.17%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.867.66'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.42%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.002'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.13%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.381'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.88%  '

$ws.Range("E17").Value = '  +5.13%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '93.13'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.27%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06602'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.48%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.48'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.60%  '

$ws.Range("E21").Value = '  -0.17%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.162'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.43%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.206.31'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.61%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.33'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.87%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.277'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.94%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.567'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +6.38%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.081.36'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.47%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '21.13'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.55%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '157.54'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.49%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '127.32'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.21%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1054'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.15%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.060'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.72%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.636'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.90%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.622'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.70%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.06757'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.93%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.519'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.91%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02397'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.82%  '

$ws.Range("E38").Value = '  +0.95%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '11.50'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.22%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.992'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.01%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6295'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.29%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.174'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.38%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.001'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.15%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.51'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.91%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6025'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.64%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.654'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.38%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.259'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.97%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '123.81'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.09%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.977'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.09%  '

$ws.Range("E50").Value = '  +0.79%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06840'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.79%  '
